$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (ID column). This shifts the existing
# A:E data to B:F and carries along per-cell styles (so the phone
# column keeps its "@" text number format, etc).
$ws.Columns.Item(1).Insert()

# --- New ID column (row numbers) ---
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3

# --- Row 1 (Yaroslav) ---
$ws.Range("C1").Value = "Shkvarlas"
$ws.Range("E1").Value = "9929845"

# --- Row 2 (Petro) ---
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = "123456"

# --- Row 3 (Ostap) ---
$ws.Range("C3").Value = "Ostapenkos"
$ws.Range("E3").Value = "2223335"

# --- Column widths ---
# New ID column ~matches the 12.63-character width used by the rest
# of the sheet, and the (now) last column (email) is widened to 25.33.
$ws.Columns.Item(1).ColumnWidth = 11.833333333333332
$ws.Columns.Item(6).ColumnWidth = 24.5

# --- Selection ---
$null = $ws.Range("C3").Select()
